# Update the multiplication problems in the practice worksheet table.
# Each old expression is unique within the document, so a simple
# Find/Replace (ReplaceAll) for each pair is safe and unambiguous.

$d = $word.ActiveDocument

$pairs = @(
    @("295×5=", "156×8="),
    @("440×8=", "628×8="),
    @("819×6=", "352×7="),
    @("300×6=", "226×2="),
    @("998×3=", "294×2="),
    @("941×4=", "418×7="),
    @("673×2=", "557×9="),
    @("105×3=", "503×6="),
    @("875×4=", "445×6="),
    @("777×6=", "497×7="),
    @("582×2=", "844×7="),
    @("213×2=", "751×9="),
    @("966×9=", "300×9="),
    @("378×9=", "585×4="),
    @("756×3=", "162×5="),
    @("380×9=", "446×4="),
    @("925×4=", "581×8="),
    @("468×2=", "950×9="),
    @("474×8=", "399×3="),
    @("703×7=", "487×4="),
    @("449×9=", "754×4="),
    @("354×2=", "555×7="),
    @("821×3=", "763×5="),
    @("804×8=", "933×3="),
    @("208×5=", "751×4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
